$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1866197183098592
$ws.Range("C2").Value = 0.5528169014084507
$ws.Range("J2").Value = 0.01056338028169014
$ws.Range("P2").Value = 0.147887323943662
$ws.Range("S2").Value = 0.102112676056338
$ws.Range("B3").Value = 0.02531645569620253
$ws.Range("C3").Value = 0.02531645569620253
$ws.Range("J3").Value = 0.0379746835443038
$ws.Range("P3").Value = 0.7025316455696202
$ws.Range("S3").Value = 0.2088607594936709
$ws.Range("P4").Value = 0.7407407407407407
$ws.Range("S4").Value = 0.2592592592592592
$ws.Range("B6").Value = 0.07011070110701106
$ws.Range("D6").Value = 0.01476014760147601
$ws.Range("F6").Value = 0.05535055350553506
$ws.Range("J6").Value = 0.2546125461254612
$ws.Range("O6").Value = 0.02214022140221402
$ws.Range("Q6").Value = 0.1328413284132841
$ws.Range("R6").Value = 0.06273062730627306
$ws.Range("S6").Value = 0.3874538745387454
$ws.Range("B7").Value = 0.05286343612334802
$ws.Range("D7").Value = 0.04405286343612335
$ws.Range("F7").Value = 0.0881057268722467
$ws.Range("J7").Value = 0.1101321585903084
$ws.Range("O7").Value = 0.02202643171806168
$ws.Range("Q7").Value = 0.1629955947136564
$ws.Range("R7").Value = 0.06607929515418502
$ws.Range("S7").Value = 0.4537444933920705
$ws.Range("B8").Value = 0.07392197125256673
$ws.Range("D8").Value = 0.01642710472279261
$ws.Range("F8").Value = 0.06776180698151951
$ws.Range("J8").Value = 0.1273100616016427
$ws.Range("O8").Value = 0.01026694045174538
$ws.Range("Q8").Value = 0.1827515400410678
$ws.Range("R8").Value = 0.1108829568788501
$ws.Range("S8").Value = 0.4106776180698152
$ws.Range("B9").Value = 0.1095238095238095
$ws.Range("D9").Value = 0.004761904761904762
$ws.Range("F9").Value = 0.05714285714285714
$ws.Range("J9").Value = 0.08571428571428572
$ws.Range("O9").Value = 0.02857142857142857
$ws.Range("Q9").Value = 0.1666666666666667
$ws.Range("R9").Value = 0.119047619047619
$ws.Range("S9").Value = 0.4285714285714285
$ws.Range("B10").Value = 0.1048
$ws.Range("D10").Value = 0.028
$ws.Range("E10").Value = 0.0008
$ws.Range("F10").Value = 0.0728
$ws.Range("J10").Value = 0.0992
$ws.Range("O10").Value = 0.0208
$ws.Range("Q10").Value = 0.2
$ws.Range("R10").Value = 0.0848
$ws.Range("S10").Value = 0.3888
$ws.Range("G11").Value = 0.1690140845070423
$ws.Range("J11").Value = 0.05633802816901409
$ws.Range("K11").Value = 0.1795774647887324
$ws.Range("L11").Value = 0.5915492957746479
$ws.Range("S11").Value = 0.00352112676056338
$ws.Range("G12").Value = 0.7724867724867724
$ws.Range("J12").Value = 0.1322751322751323
$ws.Range("K12").Value = 0.01587301587301587
$ws.Range("L12").Value = 0.06349206349206349
$ws.Range("S12").Value = 0.01587301587301587
$ws.Range("G13").Value = 0.7307692307692307
$ws.Range("J13").Value = 0.2115384615384615
$ws.Range("S13").Value = 0.0576923076923077
$ws.Range("F15").Value = 0.04807692307692308
$ws.Range("H15").Value = 0.1490384615384615
$ws.Range("I15").Value = 0.1009615384615385
$ws.Range("J15").Value = 0.3125
$ws.Range("K15").Value = 0.08653846153846154
$ws.Range("M15").Value = 0.01923076923076923
$ws.Range("N15").Value = 0.004807692307692308
$ws.Range("O15").Value = 0.0576923076923077
$ws.Range("S15").Value = 0.2211538461538461
$ws.Range("F16").Value = 0.01058201058201058
$ws.Range("H16").Value = 0.1693121693121693
$ws.Range("I16").Value = 0.06878306878306878
$ws.Range("J16").Value = 0.4603174603174603
$ws.Range("K16").Value = 0.08465608465608465
$ws.Range("M16").Value = 0.02116402116402116
$ws.Range("O16").Value = 0.03703703703703703
$ws.Range("S16").Value = 0.1481481481481481
$ws.Range("F17").Value = 0.02934537246049661
$ws.Range("H17").Value = 0.2009029345372461
$ws.Range("I17").Value = 0.09706546275395034
$ws.Range("J17").Value = 0.4198645598194131
$ws.Range("K17").Value = 0.07900677200902935
$ws.Range("M17").Value = 0.01805869074492099
$ws.Range("O17").Value = 0.03611738148984198
$ws.Range("S17").Value = 0.1196388261851016
$ws.Range("F18").Value = 0.01382488479262673
$ws.Range("H18").Value = 0.1935483870967742
$ws.Range("I18").Value = 0.08755760368663594
$ws.Range("J18").Value = 0.4147465437788018
$ws.Range("K18").Value = 0.07834101382488479
$ws.Range("M18").Value = 0.02764976958525346
$ws.Range("O18").Value = 0.07834101382488479
$ws.Range("S18").Value = 0.1059907834101382
$ws.Range("F19").Value = 0.02949852507374631
$ws.Range("H19").Value = 0.21976401179941
$ws.Range("I19").Value = 0.084070796460177
$ws.Range("J19").Value = 0.3554572271386431
$ws.Range("K19").Value = 0.1010324483775811
$ws.Range("M19").Value = 0.02286135693215339
$ws.Range("N19").Value = 0.002949852507374631
$ws.Range("O19").Value = 0.06415929203539823
$ws.Range("S19").Value = 0.1202064896755162
